$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3, shifting existing rows 3..55 down to 4..56.
$ws.Rows(3).Insert()

# Populate the newly inserted row 3 with the new weekly record.
$ws.Cells.Item(3, 1).Value = 8
$ws.Cells.Item(3, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(3, 3).Value = "Coquimbo"
$ws.Cells.Item(3, 4).Value = 45190
$ws.Cells.Item(3, 5).Value = 4
$ws.Cells.Item(3, 6).Value = 100112026
$ws.Cells.Item(3, 7).Value = "Haba"
$ws.Cells.Item(3, 8).Value = "Sin especificar"
$ws.Cells.Item(3, 9).Value = "Primera"
$ws.Cells.Item(3, 10).Value = 400
$ws.Cells.Item(3, 11).Value = 9000
$ws.Cells.Item(3, 12).Value = 10000
$ws.Cells.Item(3, 13).Value = 9500
$ws.Cells.Item(3, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(3, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(3, 16).Value = 380
$ws.Cells.Item(3, 17).Value = 25
$ws.Cells.Item(3, 18).Value = "Hortaliza"
